# Update scripts with new TPM (transcripts-per-million) values.
# The underlying NATMI ligand-receptor analysis was rerun, which:
#   1. Changed several "Target cluster" (column D) assignments for rows
#      where the original target was "Inflammatory-Mac" (now "MuSCs") or
#      "MuSCs" (now "Neutrophils").
#   2. Refreshed the computed expression / specificity statistics in
#      columns E:T across (almost) every data row to reflect the new TPM
#      inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2
    $ws.Range("I2").Value = 0.5450347805088984
    $ws.Range("J2").Value = 0.6362259982609142
    $ws.Range("M2").Value = 0.058761
    $ws.Range("N2").Value = 0.117522
    $ws.Range("O2").Value = 0.02176148226403952
    $ws.Range("P2").Value = 0.0163250711987138
    $ws.Range("Q2").Value = 0.06854703735300001
    $ws.Range("R2").Value = 0.411282224118
    $ws.Range("S2").Value = 0.01186076470932907
    $ws.Range("T2").Value = 0.01038643472008219

    # Row 3
    $ws.Range("I3").Value = 0.5450347805088984
    $ws.Range("J3").Value = 0.6362259982609142
    $ws.Range("M3").Value = 1.729918666666666
    $ws.Range("N3").Value = 5.189755999999999
    $ws.Range("O3").Value = 0.6406561219669091
    $ws.Range("P3").Value = 0.7209129882400922
    $ws.Range("Q3").Value = 2.018018744773777
    $ws.Range("S3").Value = 0.3491798688179164
    $ws.Range("T3").Value = 0.4586635856023114

    # Row 4
    $ws.Range("D4").Value = "MuSCs"
    $ws.Range("I4").Value = 0.5450347805088984
    $ws.Range("J4").Value = 0.6362259982609142
    $ws.Range("L4").Value = 1
    $ws.Range("M4").Value = 0.843062
    $ws.Range("N4").Value = 1.686124
    $ws.Range("O4").Value = 0.3122186273291074
    $ws.Range("P4").Value = 0.2342207786615281
    $ws.Range("Q4").Value = 0.9834652644593334
    $ws.Range("R4").Value = 5.900791586756
    $ws.Range("S4").Value = 0.1701700110171096
    $ws.Range("T4").Value = 0.1490173487173793

    # Row 5
    $ws.Range("D5").Value = "Neutrophils"
    $ws.Range("I5").Value = 0.5450347805088984
    $ws.Range("J5").Value = 0.6362259982609142
    $ws.Range("K5").Value = 1
    $ws.Range("L5").Value = 0.3333333333333333
    $ws.Range("M5").Value = 0.04872966666666667
    $ws.Range("N5").Value = 0.146189
    $ws.Range("O5").Value = 0.018046489625759
    $ws.Range("P5").Value = 0.02030722616589891
    $ws.Range("Q5").Value = 0.05684508911011112
    $ws.Range("R5").Value = 0.511605801991
    $ws.Range("S5").Value = 0.009835964512131666
    $ws.Range("T5").Value = 0.01291998523930919

    # Row 6
    $ws.Range("I6").Value = 0.5450347805088984
    $ws.Range("J6").Value = 0.6362259982609142
    $ws.Range("M6").Value = 0.01975833333333333
    $ws.Range("N6").Value = 0.059275
    $ws.Range("O6").Value = 0.007317278814184819
    $ws.Range("P6").Value = 0.008233935733766958
    $ws.Range("Q6").Value = 0.02304887958055556
    $ws.Range("R6").Value = 0.207439916225
    $ws.Range("S6").Value = 0.003988171452411635
    $ws.Range("T6").Value = 0.005238643981832096

    # Row 7
    $ws.Range("F7").Value = 0.3333333333333333
    $ws.Range("G7").Value = 0.053445
    $ws.Range("H7").Value = 0.160335
    $ws.Range("I7").Value = 0.02497076154086894
    $ws.Range("J7").Value = 0.02914868602301098
    $ws.Range("M7").Value = 0.058761
    $ws.Range("N7").Value = 0.117522
    $ws.Range("O7").Value = 0.02176148226403952
    $ws.Range("P7").Value = 0.0163250711987138
    $ws.Range("Q7").Value = 0.003140481645
    $ws.Range("R7").Value = 0.01884288987
    $ws.Range("S7").Value = 0.0005434007843911795
    $ws.Range("T7").Value = 0.0004758543746746081

    # Row 8
    $ws.Range("F8").Value = 0.3333333333333333
    $ws.Range("G8").Value = 0.053445
    $ws.Range("H8").Value = 0.160335
    $ws.Range("I8").Value = 0.02497076154086894
    $ws.Range("J8").Value = 0.02914868602301098
    $ws.Range("M8").Value = 1.729918666666666
    $ws.Range("N8").Value = 5.189755999999999
    $ws.Range("O8").Value = 0.6406561219669091
    $ws.Range("P8").Value = 0.7209129882400922
    $ws.Range("Q8").Value = 0.09245550313999998
    $ws.Range("R8").Value = 0.8320995282599999
    $ws.Range("S8").Value = 0.01599767125133354
    $ws.Range("T8").Value = 0.02101366634412105

    # Row 9
    $ws.Range("D9").Value = "MuSCs"
    $ws.Range("F9").Value = 0.3333333333333333
    $ws.Range("G9").Value = 0.053445
    $ws.Range("H9").Value = 0.160335
    $ws.Range("I9").Value = 0.02497076154086894
    $ws.Range("J9").Value = 0.02914868602301098
    $ws.Range("L9").Value = 1
    $ws.Range("M9").Value = 0.843062
    $ws.Range("N9").Value = 1.686124
    $ws.Range("O9").Value = 0.3122186273291074
    $ws.Range("P9").Value = 0.2342207786615281
    $ws.Range("Q9").Value = 0.04505744859
    $ws.Range("R9").Value = 0.27034469154
    $ws.Range("S9").Value = 0.007796336891652567
    $ws.Range("T9").Value = 0.006827227937270032

    # Row 10
    $ws.Range("D10").Value = "Neutrophils"
    $ws.Range("F10").Value = 0.3333333333333333
    $ws.Range("G10").Value = 0.053445
    $ws.Range("H10").Value = 0.160335
    $ws.Range("I10").Value = 0.02497076154086894
    $ws.Range("J10").Value = 0.02914868602301098
    $ws.Range("K10").Value = 1
    $ws.Range("L10").Value = 0.3333333333333333
    $ws.Range("M10").Value = 0.04872966666666667
    $ws.Range("N10").Value = 0.146189
    $ws.Range("O10").Value = 0.018046489625759
    $ws.Range("P10").Value = 0.02030722616589891
    $ws.Range("Q10").Value = 0.002604357035
    $ws.Range("R10").Value = 0.023439213315
    $ws.Range("S10").Value = 0.000450634589094593
    $ws.Range("T10").Value = 0.0005919289595080605

    # Row 11
    $ws.Range("F11").Value = 0.3333333333333333
    $ws.Range("G11").Value = 0.053445
    $ws.Range("H11").Value = 0.160335
    $ws.Range("I11").Value = 0.02497076154086894
    $ws.Range("J11").Value = 0.02914868602301098
    $ws.Range("M11").Value = 0.01975833333333333
    $ws.Range("N11").Value = 0.059275
    $ws.Range("O11").Value = 0.007317278814184819
    $ws.Range("P11").Value = 0.008233935733766958
    $ws.Range("Q11").Value = 0.001055984125
    $ws.Range("R11").Value = 0.009503857125000001
    $ws.Range("S11").Value = 0.0001827180243970613
    $ws.Range("T11").Value = 0.0002400084074372236

    # Row 12
    $ws.Range("E12").Value = 1
    $ws.Range("F12").Value = 0.5
    $ws.Range("G12").Value = 0.9203185
    $ws.Range("H12").Value = 1.840637
    $ws.Range("I12").Value = 0.4299944579502328
    $ws.Range("J12").Value = 0.3346253157160749
    $ws.Range("M12").Value = 0.058761
    $ws.Range("N12").Value = 0.117522
    $ws.Range("O12").Value = 0.02176148226403952
    $ws.Range("P12").Value = 0.0163250711987138
    $ws.Range("Q12").Value = 0.0540788353785
    $ws.Range("R12").Value = 0.216315341514
    $ws.Range("S12").Value = 0.009357316770319276
    $ws.Range("T12").Value = 0.005462782103957006

    # Row 13
    $ws.Range("E13").Value = 1
    $ws.Range("F13").Value = 0.5
    $ws.Range("G13").Value = 0.9203185
    $ws.Range("H13").Value = 1.840637
    $ws.Range("I13").Value = 0.4299944579502328
    $ws.Range("J13").Value = 0.3346253157160749
    $ws.Range("M13").Value = 1.729918666666666
    $ws.Range("N13").Value = 5.189755999999999
    $ws.Range("O13").Value = 0.6406561219669091
    $ws.Range("P13").Value = 0.7209129882400922
    $ws.Range("Q13").Value = 1.592076152428666
    $ws.Range("R13").Value = 9.552456914571998
    $ws.Range("S13").Value = 0.2754785818976593
    $ws.Range("T13").Value = 0.2412357362936599

    # Row 14
    $ws.Range("D14").Value = "MuSCs"
    $ws.Range("E14").Value = 1
    $ws.Range("F14").Value = 0.5
    $ws.Range("G14").Value = 0.9203185
    $ws.Range("H14").Value = 1.840637
    $ws.Range("I14").Value = 0.4299944579502328
    $ws.Range("J14").Value = 0.3346253157160749
    $ws.Range("L14").Value = 1
    $ws.Range("M14").Value = 0.843062
    $ws.Range("N14").Value = 1.686124
    $ws.Range("O14").Value = 0.3122186273291074
    $ws.Range("P14").Value = 0.2342207786615281
    $ws.Range("Q14").Value = 0.775885555247
    $ws.Range("R14").Value = 3.103542220988
    $ws.Range("S14").Value = 0.1342522794203453
    $ws.Range("T14").Value = 0.07837620200687873

    # Row 15
    $ws.Range("D15").Value = "Neutrophils"
    $ws.Range("E15").Value = 1
    $ws.Range("F15").Value = 0.5
    $ws.Range("G15").Value = 0.9203185
    $ws.Range("H15").Value = 1.840637
    $ws.Range("I15").Value = 0.4299944579502328
    $ws.Range("J15").Value = 0.3346253157160749
    $ws.Range("K15").Value = 1
    $ws.Range("L15").Value = 0.3333333333333333
    $ws.Range("M15").Value = 0.04872966666666667
    $ws.Range("N15").Value = 0.146189
    $ws.Range("O15").Value = 0.018046489625759
    $ws.Range("P15").Value = 0.02030722616589891
    $ws.Range("Q15").Value = 0.04484681373216667
    $ws.Range("R15").Value = 0.269080882393
    $ws.Range("S15").Value = 0.007759890524532739
    $ws.Range("T15").Value = 0.00679531196708166

    # Row 16
    $ws.Range("E16").Value = 1
    $ws.Range("F16").Value = 0.5
    $ws.Range("G16").Value = 0.9203185
    $ws.Range("H16").Value = 1.840637
    $ws.Range("I16").Value = 0.4299944579502328
    $ws.Range("J16").Value = 0.3346253157160749
    $ws.Range("M16").Value = 0.01975833333333333
    $ws.Range("N16").Value = 0.059275
    $ws.Range("O16").Value = 0.007317278814184819
    $ws.Range("P16").Value = 0.008233935733766958
    $ws.Range("Q16").Value = 0.01818395969583333
    $ws.Range("R16").Value = 0.109103758175
    $ws.Range("S16").Value = 0.003146389337376123
    $ws.Range("T16").Value = 0.002755283344497639
